# handover_data.xlsx update
# - Row 2: clear ApprovalToReceive (Q2) and CompletionDate (S2), flip Status (T2) from
#   "Rejected" to "Pending" (form no longer completed / receiver hasn't actioned it yet).
# - Rows 3-7: five additional handover-form submissions appended below the existing one
#   (transfer-progress / transfer-history test rows + one real "priority test left" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: amend in place ----------------------------------------------------
$ws.Range("Q2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("T2").Value = "Pending"

# ---- New rows 3-7 --------------------------------------------------------------
# Each hashtable maps column letter -> value for that row. Missing keys are left blank.
$rows = @{
    3 = @{
        A = "asd"; B = 111111111111; C = "Drone Equipment"; D = "Name7"; E = "Make7"
        F = "Model7"; G = "asd"; H = "a"; I = "asd"; J = "Umar"; K = "Umar"
        L = "Not OK"; M = "-"; N = "Good"; O = "-"; P = "YES"
        R = "2024-06-07 12:54:38"; T = "Pending"
    }
    4 = @{
        A = "sdsdsa"; B = 111111111111; C = "Drone Equipment"; D = "Name7"; E = "Make7"
        F = "Model7"; G = "Serial7"; H = "ad"; I = "SOI TRIPURA"; J = "Umar"; K = "asd"
        L = "Not OK"; M = "-"; N = "Good"; O = "-"
        R = "2024-06-07 12:54:38"; T = "Pending"
    }
    5 = @{
        A = "das"; B = "aaaaaaaaaaaa"; C = "Drone Equipment"; D = "Name7"; E = "Make7"
        F = "Model7"; G = "Serial7"; H = "SOI ASSAM"; I = "SOI TRIPURA"; J = "sdasd"; K = "Umar"
        L = "Not OK"; M = "-"; N = "Good"; O = "-"; P = 1
        R = "2024-06-07 12:54:38"; S = "2024-06-07 12:54:38"; T = "Pending"
    }
    6 = @{
        A = "s"; B = 111111111111; C = "Drone Equipment"; D = "Name7"; E = "Make7"
        F = "Model7"; G = "Serial7"; H = "SOI ASSAM"; I = "SOI TRIPURA"; J = "Umar"; K = "asd"
        L = "Not OK"; M = "-"; N = "Good"; O = "-"
        R = "2024-06-07 12:54:38"; S = "2024-06-07 12:54:38"; T = "Pending"
    }
    7 = @{
        A = "ca1bd423"; B = "123412341234"; C = "Electronics"; D = "Name1"; E = "Make1"
        F = "Model1"; G = "Serial1"; H = "SOI ASSAM"; I = "SOI TRIPURA"; J = "Hammad"; K = "Fahad"
        L = "Good"; M = "badhiya"; N = "Not Ok"; P = 1; Q = 1
        R = "2024-06-07 18:48:17"; S = "2024-06-07 19:25:29"; T = "Approved"
    }
}

# EwayBillNo values that are pure digit strings would otherwise be auto-parsed as numbers;
# tag them so they are written back as text (matches the source export, which stores this
# particular EwayBillNo as text rather than a number).
$textColumnsByRow = @{ 7 = @("B") }

foreach ($r in ($rows.Keys | Sort-Object)) {
    $rowData = $rows[$r]
    $textCols = $textColumnsByRow[$r]
    foreach ($col in $rowData.Keys) {
        $value = $rowData[$col]
        $cell = $ws.Range("$col$r")
        if ($textCols -contains $col) {
            # Force text storage for an all-digit EwayBillNo (otherwise it would be
            # auto-parsed back into a number) and then drop the quote-prefix display
            # style the leading apostrophe implies, restoring the default cell style.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
